{"js": "// Update the date line and the 25 division-problem answers in the table.\n// Each original text value is unique in the document, so a scoped,\n// case-sensitive search-and-replace for each string is unambiguous and\n// safe to run unconditionally.\nconst replacements = [\n  [\"2025-02-12 Wednesday\", \"2025-02-13 Thursday\"],\n  [\"963\u00f74=240, 3\", \"914\u00f74=228, 2\"],\n  [\"628\u00f75=125, 3\", \"785\u00f78=98, 1\"],\n  [\"207\u00f74=51, 3\", \"451\u00f78=56, 3\"],\n  [\"296\u00f72=148, 0\", \"704\u00f74=176, 0\"],\n  [\"158\u00f77=22, 4\", \"180\u00f74=45, 0\"],\n  [\"800\u00f79=88, 8\", \"125\u00f72=62, 1\"],\n  [\"434\u00f72=217, 0\", \"748\u00f78=93, 4\"],\n  [\"460\u00f72=230, 0\", \"685\u00f78=85, 5\"],\n  [\"992\u00f73=330, 2\", \"762\u00f74=190, 2\"],\n  [\"849\u00f72=424, 1\", \"703\u00f78=87, 7\"],\n  [\"755\u00f73=251, 2\", \"398\u00f76=66, 2\"],\n  [\"862\u00f78=107, 6\", \"714\u00f74=178, 2\"],\n  [\"940\u00f74=235, 0\", \"312\u00f78=39, 0\"],\n  [\"138\u00f76=23, 0\", \"686\u00f76=114, 2\"],\n  [\"205\u00f77=29, 2\", \"452\u00f79=50, 2\"],\n  [\"986\u00f77=140, 6\", \"199\u00f75=39, 4\"],\n  [\"806\u00f76=134, 2\", \"701\u00f74=175, 1\"],\n  [\"695\u00f76=115, 5\", \"689\u00f72=344, 1\"],\n  [\"111\u00f74=27, 3\", \"920\u00f78=115, 0\"],\n  [\"315\u00f73=105, 0\", \"337\u00f76=56, 1\"],\n  [\"722\u00f75=144, 2\", \"973\u00f77=139, 0\"],\n  [\"847\u00f73=282, 1\", \"969\u00f74=242, 1\"],\n  [\"881\u00f77=125, 6\", \"292\u00f72=146, 0\"],\n  [\"639\u00f72=319, 1\", \"948\u00f72=474, 0\"],\n  [\"851\u00f73=283, 2\", \"198\u00f79=22, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division-problem answers in the table.\n# Each original text value is unique in the document, so a document-wide\n# Find/Replace (one execution per pair, Replace:=wdReplaceAll) is safe and\n# unambiguous for every entry.\n\n$pairs = @(\n    @(\"2025-02-12 Wednesday\", \"2025-02-13 Thursday\"),\n    @(\"963\u00f74=240, 3\", \"914\u00f74=228, 2\"),\n    @(\"628\u00f75=125, 3\", \"785\u00f78=98, 1\"),\n    @(\"207\u00f74=51, 3\", \"451\u00f78=56, 3\"),\n    @(\"296\u00f72=148, 0\", \"704\u00f74=176, 0\"),\n    @(\"158\u00f77=22, 4\", \"180\u00f74=45, 0\"),\n    @(\"800\u00f79=88, 8\", \"125\u00f72=62, 1\"),\n    @(\"434\u00f72=217, 0\", \"748\u00f78=93, 4\"),\n    @(\"460\u00f72=230, 0\", \"685\u00f78=85, 5\"),\n    @(\"992\u00f73=330, 2\", \"762\u00f74=190, 2\"),\n    @(\"849\u00f72=424, 1\", \"703\u00f78=87, 7\"),\n    @(\"755\u00f73=251, 2\", \"398\u00f76=66, 2\"),\n    @(\"862\u00f78=107, 6\", \"714\u00f74=178, 2\"),\n    @(\"940\u00f74=235, 0\", \"312\u00f78=39, 0\"),\n    @(\"138\u00f76=23, 0\", \"686\u00f76=114, 2\"),\n    @(\"205\u00f77=29, 2\", \"452\u00f79=50, 2\"),\n    @(\"986\u00f77=140, 6\", \"199\u00f75=39, 4\"),\n    @(\"806\u00f76=134, 2\", \"701\u00f74=175, 1\"),\n    @(\"695\u00f76=115, 5\", \"689\u00f72=344, 1\"),\n    @(\"111\u00f74=27, 3\", \"920\u00f78=115, 0\"),\n    @(\"315\u00f73=105, 0\", \"337\u00f76=56, 1\"),\n    @(\"722\u00f75=144, 2\", \"973\u00f77=139, 0\"),\n    @(\"847\u00f73=282, 1\", \"969\u00f74=242, 1\"),\n    @(\"881\u00f77=125, 6\", \"292\u00f72=146, 0\"),\n    @(\"639\u00f72=319, 1\", \"948\u00f72=474, 0\"),\n    @(\"851\u00f73=283, 2\", \"198\u00f79=22, 0\")\n)\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $searchRange = $d.Content\n    $found = $searchRange.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, $wdReplaceAll)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
